$d = $word.ActiveDocument

$replacements = @(
    @{old="572×2=1144"; new="839×9=7551"},
    @{old="386×6=2316"; new="655×8=5240"},
    @{old="545×5=2725"; new="145×7=1015"},
    @{old="517×9=4653"; new="515×7=3605"},
    @{old="703×8=5624"; new="103×6=618"},
    @{old="575×4=2300"; new="142×3=426"},
    @{old="351×8=2808"; new="988×6=5928"},
    @{old="994×2=1988"; new="125×4=500"},
    @{old="901×4=3604"; new="418×9=3762"},
    @{old="579×8=4632"; new="459×3=1377"},
    @{old="568×3=1704"; new="172×4=688"},
    @{old="604×2=1208"; new="773×8=6184"},
    @{old="562×9=5058"; new="556×3=1668"},
    @{old="341×9=3069"; new="434×7=3038"},
    @{old="151×7=1057"; new="485×6=2910"},
    @{old="276×5=1380"; new="526×9=4734"},
    @{old="158×9=1422"; new="233×8=1864"},
    @{old="190×4=760"; new="563×7=3941"},
    @{old="979×4=3916"; new="645×5=3225"},
    @{old="984×9=8856"; new="395×9=3555"},
    @{old="880×2=1760"; new="186×4=744"},
    @{old="683×8=5464"; new="206×8=1648"},
    @{old="461×7=3227"; new="631×4=2524"},
    @{old="876×7=6132"; new="839×6=5034"},
    @{old="428×8=3424"; new="991×7=6937"}
)

foreach ($r in $replacements) {
    $d.Content.Find.Execute($r.old, $true, $false, $false, $false, $false,
                             $true, 1, $false, $r.new, 2)
}
